$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.753.66"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "2.288.26"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'110.90"
$ws.Range("E5").Value = "  +14.11%  "

$ws.Range("D6").Value = "'268.54"

$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").Value = "'47.71"
$ws.Range("E10").Value = "  +4.71%  "

$ws.Range("D11").Value = "'0.0945"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("E12").Value = "  +13.94%  "

$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").Value = "'15.74"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").Value = "2.629.95"
$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").Value = "'0.846"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").Value = "2.271.15"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").Value = "43.618.97"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("D20").Value = "'6.72"
$ws.Range("E20").Value = "  +8.36%  "

$ws.Range("D21").Value = "'72.13"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "'2.43"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("D23").Value = "'9.86"
$ws.Range("E23").Value = "  +8.67%  "

$ws.Range("D24").Value = "'231.84"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").Value = "'2.76"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'11.59"
$ws.Range("E27").Value = "  +3.05%  "

$ws.Range("D28").Value = "'41.62"
$ws.Range("E28").Value = "  +5.13%  "

$ws.Range("D29").Value = "'3.40"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").Value = "'175.31"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").Value = "'21.45"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").Value = "'0.0922"
$ws.Range("E33").Value = "  +2.47%  "

$ws.Range("D34").Value = "'5.60"
$ws.Range("E34").Value = "  +3.85%  "

$ws.Range("D35").Value = "'0.127"
$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("D36").Value = "'4.67"
$ws.Range("E36").Value = "  +5.62%  "

$ws.Range("D37").Value = "'0.0362"
$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("D39").Value = "'3.83"
$ws.Range("E39").Value = "  +13.16%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.40"
$ws.Range("E40").Value = "  +2.50%  "

$ws.Range("D41").Value = "'0.241"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'73.13"
$ws.Range("E42").Value = "  +13.04%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'13.66"
$ws.Range("E43").Value = "  +10.91%  "

$ws.Range("D44").Value = "'6.31"
$ws.Range("E44").Value = "  +22.35%  "

$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").Value = "'1.37"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").Value = "'8.71"
$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").Value = "'102.49"
$ws.Range("E48").Value = "  +5.36%  "

$ws.Range("D49").Value = "'0.0994"
$ws.Range("E49").Value = "  -2.36%  "

$ws.Range("D50").Value = "'1.22"
$ws.Range("E50").Value = "  +2.14%  "

$ws.Range("D51").Value = "'0.449"
$ws.Range("E51").Value = "  +4.88%  "
